$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<a>"
$ws.Range("C2").Value = 17

# Row 3
$ws.Range("C3").Value = 24

# Row 4
$ws.Range("C4").Value = 11

# Row 5
$ws.Range("C5").Value = 15

# Row 6
$ws.Range("B6").Value = "<seven>"
$ws.Range("C6").Value = 16

# Row 8
$ws.Range("C8").Value = 19

# Row 9
$ws.Range("C9").Value = 13

# Row 10
$ws.Range("C10").Value = 9

# Row 11
$ws.Range("C11").Value = 16

# Row 12
$ws.Range("B12").Value = "<hen>"
$ws.Range("C12").Value = 12

# Row 13
$ws.Range("C13").Value = 20

# Row 14
$ws.Range("C14").Value = 9

# Row 15
$ws.Range("C15").Value = 13

# Row 16
$ws.Range("C16").Value = 11

# Row 17
$ws.Range("B17").Value = "<left>"
$ws.Range("C17").Value = 16

# Row 18
$ws.Range("C18").Value = 12
